$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update D1 header value
$ws.Range("D1").Value = "Application: Timups Website"

# 2. Update E9 remark text (Read More button remark rewritten)
$ws.Range("E9").Value = '"Read More" button hasn''t been implemented and instead it reloads the same "About Us" page'

# 3. Add new row 10 data (new test scenario about Timups logo)
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Timups logo"
$ws.Range("C10").Value = "Pass"
$ws.Range("E10").Value = "This has been tested from every single page developped and has worked from all of them "
$ws.Range("D10").Value = "When the user clicks on Timups logo, then it brings them back to the Homepage"

# Copy formatting from row 9 (same column) down to row 10 so the new row matches existing styling
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)

# Re-apply the "Pass" (green) styling used elsewhere for C column, matching C5/C6
$ws.Range("C5").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)

$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# Re-set values after paste-special in case formatting paste cleared them
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Timups logo"
$ws.Range("C10").Value = "Pass"
$ws.Range("E10").Value = "This has been tested from every single page developped and has worked from all of them "
$ws.Range("D10").Value = "When the user clicks on Timups logo, then it brings them back to the Homepage"

$ws.Rows.Item(10).RowHeight = 49.2

# 4. Update the active selection to E10
$ws.Range("E10").Select()
